# Book1.xlsx edit: add a "Q" column (HEX2DEC of the opcode's hex encoding
# already computed in column E) down through the instruction table, and
# refresh the fill-down of the existing E / G / P helper formulas so Excel
# re-groups them as shared formulas (this happens naturally when a column
# of formulas is filled down in one action).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-fill the existing helper formulas (column E, G, P) in the same
#     three 32-row blocks the author filled them in, so Excel re-derives
#     the shared-formula groupings across the whole table. ---

# Block 1: rows 2-33
$ws.Range("E2:E33").Formula = "=CONCAT(TEXT(BIN2HEX(LEFT(A2,4)), ""0""), TEXT(BIN2HEX(CONCAT(RIGHT(A2), B2, C2,D2)), ""0""))"
$ws.Range("G2:G33").Formula = "=B2"
$ws.Range("P2:P33").Formula = "=CONCAT(TEXT(BIN2HEX(CONCAT(F2,G2,H2,LEFT(I2))),""0""), TEXT(BIN2HEX(CONCAT(RIGHT(I2,3), J2)), ""0""), TEXT(BIN2HEX(CONCAT(K2,L2,M2,N2)),""0""), TEXT(BIN2HEX(CONCAT(O2,""000"")), ""0""))"

# Block 2: rows 34-65
$ws.Range("E34:E65").Formula = "=CONCAT(TEXT(BIN2HEX(LEFT(A34,4)), ""0""), TEXT(BIN2HEX(CONCAT(RIGHT(A34), B34, C34,D34)), ""0""))"
$ws.Range("G34:G65").Formula = "=B34"
$ws.Range("P34:P65").Formula = "=CONCAT(TEXT(BIN2HEX(CONCAT(F34,G34,H34,LEFT(I34))),""0""), TEXT(BIN2HEX(CONCAT(RIGHT(I34,3), J34)), ""0""), TEXT(BIN2HEX(CONCAT(K34,L34,M34,N34)),""0""), TEXT(BIN2HEX(CONCAT(O34,""000"")), ""0""))"

# Block 3: rows 66-73 (E/P filled through the 32-row block extent 66-97,
# matching the original fill action; the sheet only has data through row
# 73 so the extra rows are removed again immediately after).
$ws.Range("E66:E97").Formula = "=CONCAT(TEXT(BIN2HEX(LEFT(A66,4)), ""0""), TEXT(BIN2HEX(CONCAT(RIGHT(A66), B66, C66,D66)), ""0""))"
$ws.Range("G66:G73").Formula = "=B66"
$ws.Range("P66:P97").Formula = "=CONCAT(TEXT(BIN2HEX(CONCAT(F66,G66,H66,LEFT(I66))),""0""), TEXT(BIN2HEX(CONCAT(RIGHT(I66,3), J66)), ""0""), TEXT(BIN2HEX(CONCAT(K66,L66,M66,N66)),""0""), TEXT(BIN2HEX(CONCAT(O66,""000"")), ""0""))"
$ws.Range("A74:A97").EntireRow.Delete()

# --- New column Q: decimal value of the hex-encoded opcode in column E ---
$ws.Range("Q2").Formula = "=HEX2DEC(E2)"
$ws.Range("Q3:Q66").Formula = "=HEX2DEC(E3)"
$ws.Range("Q67:Q73").Formula = "=HEX2DEC(E67)"

# --- Restore the selection the author left the sheet in ---
$ws.Range("R71").Select() | Out-Null
